# kon/koff data analysis pass: rename the "SE" (standard error) headers to
# "SD" (standard deviation) now that the error values being reported for the
# kon/koff/Kd columns are standard deviations rather than standard errors,
# and tidy up a couple of stray cell styles left over from the original
# formatting pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: "kon SE" / "koff SE" / "Kd SE" -> "kon SD" / "koff SD" / "Kd SD"
$ws.Range("G1").Value = "kon SD"
$ws.Range("I1").Value = "koff SD"
$ws.Range("K1").Value = "Kd SD"

# A18 and C19 still carried the old border/vertical-alignment style left
# over from an earlier pass (the border row that used to separate sections).
# Clear them back to the plain/default style used by the rest of the table.
$ws.Range("A18").ClearFormats()
$ws.Range("C19").ClearFormats()

# Leave the selection where work left off.
$ws.Range("B9").Select() | Out-Null
